# Auto-generated edit script applying the Odin_Profits.xlsx data-refresh diff
# (market price / profit recompute across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 79.59999999999999
$ws.Range("I6").Value = 79.59999999999999
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 238.8
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -126.8
$ws.Range("H15").Value = 405364.28
$ws.Range("I15").Value = 405364.28
$ws.Range("K15").Value = 1216092.84
$ws.Range("M15").Value = -1215923.84
$ws.Range("H19").Value = 3158.125
$ws.Range("I19").Value = 888.6667
$ws.Range("K19").Value = 888.6667
$ws.Range("M19").Value = -713.6667
$ws.Range("H33").Value = 693.1667
$ws.Range("I33").Value = 636.9
$ws.Range("K33").Value = 636.9
$ws.Range("M33").Value = -407.9
$ws.Range("H38").Value = 125.25
$ws.Range("I38").Value = 125.25
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 375.75
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -3.75
$ws.Range("H40").Value = 4161.6665
$ws.Range("I40").Value = 3995
$ws.Range("J40").Value = 4245
$ws.Range("K40").Value = 3995
$ws.Range("L40").Value = 4245
$ws.Range("M40").Value = -3820
$ws.Range("N40").Value = -4595
$ws.Range("H41").Value = 522
$ws.Range("I41").Value = 650
$ws.Range("J41").Value = 479.33334
$ws.Range("K41").Value = 650
$ws.Range("L41").Value = 479.33334
$ws.Range("M41").Value = -210
$ws.Range("N41").Value = -1359.33334
$ws.Range("H43").Value = 6562.143
$ws.Range("I43").Value = 2000
$ws.Range("J43").Value = 6913.077
$ws.Range("K43").Value = 2000
$ws.Range("L43").Value = 6913.077
$ws.Range("M43").Value = -1931
$ws.Range("N43").Value = -7051.077
$ws.Range("H51").Value = 41096
$ws.Range("J51").Value = 6649
$ws.Range("L51").Value = 6649
$ws.Range("N51").Value = -7617
$ws.Range("H55").Value = 350.45456
$ws.Range("I55").Value = 179.2
$ws.Range("J55").Value = 493.16666
$ws.Range("K55").Value = 179.2
$ws.Range("L55").Value = 493.16666
$ws.Range("M55").Value = 34.80000000000001
$ws.Range("N55").Value = -921.16666
$ws.Range("H64").Value = 55559428
$ws.Range("I64").Value = 55559428
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 55559428
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -55559180
$ws.Range("H67").Value = 55559428
$ws.Range("I67").Value = 55559428
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 55559428
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -55558570
$ws.Range("H70").Value = 1336.875
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1336.875
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 4010.625
$ws.Range("N70").Value = -4550.625
$ws.Range("H73").Value = 1336.875
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1336.875
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 4010.625
$ws.Range("N73").Value = -5882.625
$ws.Range("H76").Value = 76929710
$ws.Range("I76").Value = 125005820
$ws.Range("K76").Value = 125005820
$ws.Range("M76").Value = -125005505
$ws.Range("H79").Value = 76929710
$ws.Range("I79").Value = 125005820
$ws.Range("K79").Value = 125005820
$ws.Range("M79").Value = -125004728
$ws.Range("H80").Value = 843.38464
$ws.Range("I80").Value = 351.75
$ws.Range("J80").Value = 1630
$ws.Range("K80").Value = 1055.25
$ws.Range("L80").Value = 4890
$ws.Range("M80").Value = -57.25
$ws.Range("N80").Value = -6886
$ws.Range("H83").Value = 843.38464
$ws.Range("I83").Value = 351.75
$ws.Range("J83").Value = 1630
$ws.Range("K83").Value = 3165.75
$ws.Range("L83").Value = 14670
$ws.Range("M83").Value = 1826.25
$ws.Range("N83").Value = -24654
$ws.Range("H86").Value = 100001950
$ws.Range("I86").Value = 200001550
$ws.Range("J86").Value = 2358.6
$ws.Range("K86").Value = 200001550
$ws.Range("L86").Value = 2358.6
$ws.Range("M86").Value = -200000427
$ws.Range("N86").Value = -4604.6
$ws.Range("H89").Value = 100001950
$ws.Range("I89").Value = 200001550
$ws.Range("J89").Value = 2358.6
$ws.Range("K89").Value = 1000007750
$ws.Range("L89").Value = 11793
$ws.Range("M89").Value = -1000002134
$ws.Range("N89").Value = -23025
$ws.Range("H100").Value = 1630.2667
$ws.Range("I100").Value = 987.1818
$ws.Range("K100").Value = 987.1818
$ws.Range("M100").Value = -446.1818
$ws.Range("H103").Value = 1608.625
$ws.Range("I103").Value = 733.6
$ws.Range("J103").Value = 3067
$ws.Range("K103").Value = 2200.8
$ws.Range("L103").Value = 9201
$ws.Range("M103").Value = -1614.8
$ws.Range("N103").Value = -10373
$ws.Range("H111").Value = 1245.375
$ws.Range("I111").Value = 1077.1666
$ws.Range("J111").Value = 1750
$ws.Range("K111").Value = 3231.4998
$ws.Range("L111").Value = 5250
$ws.Range("M111").Value = -164.4998000000001
$ws.Range("N111").Value = -11384
$ws.Range("H112").Value = 3201.85
$ws.Range("J112").Value = 2370.3684
$ws.Range("L112").Value = 7111.1052
$ws.Range("N112").Value = -9327.1052
$ws.Range("H116").Value = 8563317
$ws.Range("I116").Value = 10116648
$ws.Range("J116").Value = 19998
$ws.Range("K116").Value = 10116648
$ws.Range("L116").Value = 19998
$ws.Range("M116").Value = -10113206
$ws.Range("N116").Value = -26882
$ws.Range("H127").Value = 6976.25
$ws.Range("I127").Value = 7294.8887
$ws.Range("J127").Value = 4108.5
$ws.Range("K127").Value = 21884.6661
$ws.Range("L127").Value = 12325.5
$ws.Range("M127").Value = -16924.6661
$ws.Range("N127").Value = -22245.5
$ws.Range("H132").Value = 360818.94
$ws.Range("I132").Value = 429969.44
$ws.Range("J132").Value = 15066.5
$ws.Range("K132").Value = 1289908.32
$ws.Range("L132").Value = 45199.5
$ws.Range("M132").Value = -1287378.32
$ws.Range("N132").Value = -50259.5
$ws.Range("H135").Value = 3628.55
$ws.Range("I135").Value = 2137.8462
$ws.Range("K135").Value = 19240.6158
$ws.Range("M135").Value = -16705.6158
$ws.Range("N6").ClearContents()
$ws.Range("N38").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("N67").ClearContents()
$ws.Range("M70").ClearContents()
$ws.Range("M73").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4678.6875
$ws.Range("I2").Value = 4289.72
$ws.Range("J2").Value = 6067.857
$ws.Range("K2").Value = 4289.72
$ws.Range("L2").Value = 6067.857
$ws.Range("M2").Value = -4176.72
$ws.Range("N2").Value = -6293.857
$ws.Range("H5").Value = 824.5
$ws.Range("J5").Value = 750
$ws.Range("L5").Value = 750
$ws.Range("N5").Value = -974
$ws.Range("H32").Value = 2098610.2
$ws.Range("I32").Value = 816.8333
$ws.Range("K32").Value = 816.8333
$ws.Range("M32").Value = -529.8333
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("H45").Value = 1691.3823
$ws.Range("I45").Value = 1690.1538
$ws.Range("J45").Value = 1695.375
$ws.Range("K45").Value = 1690.1538
$ws.Range("L45").Value = 1695.375
$ws.Range("M45").Value = -1313.1538
$ws.Range("N45").Value = -2449.375
$ws.Range("H61").Value = 4211.364
$ws.Range("I61").Value = 4227.778
$ws.Range("K61").Value = 4227.778
$ws.Range("M61").Value = -4015.778
$ws.Range("H63").Value = 3733.35
$ws.Range("J63").Value = 4792.727
$ws.Range("L63").Value = 4792.727
$ws.Range("N63").Value = -6164.727
$ws.Range("H66").Value = 3733.35
$ws.Range("J66").Value = 4792.727
$ws.Range("L66").Value = 23963.635
$ws.Range("N66").Value = -30827.635
$ws.Range("H74").Value = 2156.0557
$ws.Range("I74").Value = 1301.1
$ws.Range("K74").Value = 1301.1
$ws.Range("M74").Value = -427.0999999999999
$ws.Range("H77").Value = 2156.0557
$ws.Range("I77").Value = 1301.1
$ws.Range("K77").Value = 6505.5
$ws.Range("M77").Value = -2137.5
$ws.Range("H97").Value = 762.3333
$ws.Range("J97").Value = 756.9091
$ws.Range("L97").Value = 756.9091
$ws.Range("N97").Value = -1748.9091
$ws.Range("H102").Value = 10788.4
$ws.Range("I102").Value = 9735.5
$ws.Range("K102").Value = 9735.5
$ws.Range("M102").Value = -8113.5
$ws.Range("H116").Value = 4678.6875
$ws.Range("I116").Value = 4289.72
$ws.Range("J116").Value = 6067.857
$ws.Range("K116").Value = 4289.72
$ws.Range("L116").Value = 6067.857
$ws.Range("M116").Value = -1995.72
$ws.Range("N116").Value = -10655.857
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H136").Value = 4211.364
$ws.Range("I136").Value = 4227.778
$ws.Range("K136").Value = 12683.334
$ws.Range("M136").Value = -10133.334
$ws.Range("H137").Value = 99994
$ws.Range("J137").Value = 99994
$ws.Range("L137").Value = 99994
$ws.Range("N137").Value = -110194
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
$ws.Range("N135").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4678.6875
$ws.Range("I3").Value = 4289.72
$ws.Range("J3").Value = 6067.857
$ws.Range("K3").Value = 4289.72
$ws.Range("L3").Value = 6067.857
$ws.Range("M3").Value = -4175.72
$ws.Range("N3").Value = -6295.857
$ws.Range("H4").Value = 824.5
$ws.Range("J4").Value = 750
$ws.Range("L4").Value = 750
$ws.Range("N4").Value = -980
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("H63").Value = 74999.91
$ws.Range("J63").Value = 74999.91
$ws.Range("L63").Value = 74999.91
$ws.Range("N63").Value = -76371.91
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("H66").Value = 74999.91
$ws.Range("J66").Value = 74999.91
$ws.Range("L66").Value = 224999.73
$ws.Range("N66").Value = -231863.73
$ws.Range("H82").Value = 48333.332
$ws.Range("J82").Value = 45000
$ws.Range("L82").Value = 45000
$ws.Range("N82").Value = -45766
$ws.Range("H85").Value = 48333.332
$ws.Range("J85").Value = 45000
$ws.Range("L85").Value = 45000
$ws.Range("N85").Value = -47652
$ws.Range("H94").Value = 742.4483
$ws.Range("I94").Value = 831.7273
$ws.Range("K94").Value = 831.7273
$ws.Range("M94").Value = -380.7273
$ws.Range("H99").Value = 7329.881
$ws.Range("I99").Value = 8171
$ws.Range("K99").Value = 8171
$ws.Range("M99").Value = -6673
$ws.Range("H105").Value = 3589.762
$ws.Range("I105").Value = 4512.5835
$ws.Range("K105").Value = 4512.5835
$ws.Range("M105").Value = -2765.5835
$ws.Range("H134").Value = 3405560.2
$ws.Range("I134").Value = 3665237.8
$ws.Range("J134").Value = 29750
$ws.Range("K134").Value = 10995713.4
$ws.Range("L134").Value = 89250
$ws.Range("M134").Value = -10993178.4
$ws.Range("N134").Value = -94320
$ws.Range("H138").Value = 84999
$ws.Range("J138").Value = 84999
$ws.Range("L138").Value = 84999
$ws.Range("N138").Value = -95279
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 55562470
$ws.Range("I16").Value = 125004060
$ws.Range("J16").Value = 9199.6
$ws.Range("K16").Value = 125004060
$ws.Range("L16").Value = 9199.6
$ws.Range("M16").Value = -125003773
$ws.Range("N16").Value = -9773.6
$ws.Range("H22").Value = 2981246.8
$ws.Range("I22").Value = 4465870
$ws.Range("K22").Value = 4465870
$ws.Range("M22").Value = -4465520
$ws.Range("H25").Value = 55.5
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("H31").Value = 2457.1538
$ws.Range("I31").Value = 1325.5
$ws.Range("K31").Value = 1325.5
$ws.Range("M31").Value = -1030.5
$ws.Range("H34").Value = 2457.1538
$ws.Range("I34").Value = 1325.5
$ws.Range("K34").Value = 1325.5
$ws.Range("M34").Value = -1123.5
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("H58").Value = 58831510
$ws.Range("I58").Value = 100005030
$ws.Range("J58").Value = 12192.571
$ws.Range("K58").Value = 100005030
$ws.Range("L58").Value = 12192.571
$ws.Range("M58").Value = -100004827
$ws.Range("N58").Value = -12598.571
$ws.Range("H107").Value = 699.8276
$ws.Range("I107").Value = 537.3333
$ws.Range("K107").Value = 537.3333
$ws.Range("M107").Value = 1382.6667
$ws.Range("H113").Value = 55562470
$ws.Range("I113").Value = 125004060
$ws.Range("J113").Value = 9199.6
$ws.Range("K113").Value = 125004060
$ws.Range("L113").Value = 9199.6
$ws.Range("M113").Value = -125001890
$ws.Range("N113").Value = -13539.6
$ws.Range("H132").Value = 6402.2964
$ws.Range("I132").Value = 5423.35
$ws.Range("J132").Value = 9199.286
$ws.Range("K132").Value = 16270.05
$ws.Range("L132").Value = 27597.858
$ws.Range("M132").Value = -13740.05
$ws.Range("N132").Value = -32657.858
$ws.Range("H136").Value = 58831510
$ws.Range("I136").Value = 100005030
$ws.Range("J136").Value = 12192.571
$ws.Range("K136").Value = 300015090
$ws.Range("L136").Value = 36577.713
$ws.Range("M136").Value = -300012540
$ws.Range("N136").Value = -41677.713
$ws.Range("H141").Value = 154141.28
$ws.Range("J141").Value = 49797.8
$ws.Range("L141").Value = 49797.8
$ws.Range("N141").Value = -60157.8
$ws.Range("N25").ClearContents()
$ws.Range("M44").ClearContents()
$ws.Range("N44").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 186.4375
$ws.Range("J23").Value = 532.5
$ws.Range("L23").Value = 1597.5
$ws.Range("N23").Value = -2067.5
$ws.Range("H26").Value = 331.47058
$ws.Range("I26").Value = 339.75
$ws.Range("K26").Value = 1019.25
$ws.Range("M26").Value = -731.25
$ws.Range("H34").Value = 411.83334
$ws.Range("J34").Value = 887
$ws.Range("L34").Value = 2661
$ws.Range("N34").Value = -2829
$ws.Range("H38").Value = 47.166668
$ws.Range("I38").Value = 38.25
$ws.Range("K38").Value = 114.75
$ws.Range("M38").Value = 232.25
$ws.Range("H56").Value = 7878.1177
$ws.Range("I56").Value = 7878.1177
$ws.Range("K56").Value = 7878.1177
$ws.Range("M56").Value = -7348.1177
$ws.Range("H92").Value = 266.57693
$ws.Range("J92").Value = 470.2
$ws.Range("L92").Value = 1410.6
$ws.Range("N92").Value = -3906.6
$ws.Range("H98").Value = 1440.6364
$ws.Range("I98").Value = 2295.3333
$ws.Range("K98").Value = 6885.999899999999
$ws.Range("M98").Value = -5387.999899999999
$ws.Range("H122").Value = 218294.78
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 218294.78
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 1964653.02
$ws.Range("N122").Value = -1969553.02
$ws.Range("M122").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 100016390
$ws.Range("J80").Value = 26791.5
$ws.Range("L80").Value = 26791.5
$ws.Range("N80").Value = -28787.5
$ws.Range("H83").Value = 100016390
$ws.Range("J83").Value = 26791.5
$ws.Range("L83").Value = 133957.5
$ws.Range("N83").Value = -143941.5
$ws.Range("H102").Value = 4496.2085
$ws.Range("I102").Value = 4070.5264
$ws.Range("K102").Value = 4070.5264
$ws.Range("M102").Value = -2448.5264
$ws.Range("H122").Value = 9220.916999999999
$ws.Range("I122").Value = 5001.6
$ws.Range("J122").Value = 12234.714
$ws.Range("K122").Value = 15004.8
$ws.Range("L122").Value = 36704.142
$ws.Range("M122").Value = -12554.8
$ws.Range("N122").Value = -41604.142
$ws.Range("H132").Value = 111114900
$ws.Range("I132").Value = 250002340
$ws.Range("J132").Value = 4944.6
$ws.Range("K132").Value = 750007020
$ws.Range("L132").Value = 14833.8
$ws.Range("M132").Value = -750004490
$ws.Range("N132").Value = -19893.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4616.931
$ws.Range("I7").Value = 3173.75
$ws.Range("J7").Value = 5635.647
$ws.Range("K7").Value = 3173.75
$ws.Range("L7").Value = 5635.647
$ws.Range("M7").Value = -3061.75
$ws.Range("N7").Value = -5859.647
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("H46").Value = 17858616
$ws.Range("I46").Value = 802.2
$ws.Range("J46").Value = 38463784
$ws.Range("K46").Value = 802.2
$ws.Range("L46").Value = 38463784
$ws.Range("M46").Value = -614.2
$ws.Range("N46").Value = -38464160
$ws.Range("H55").Value = 1880.1471
$ws.Range("I55").Value = 1185.4286
$ws.Range("J55").Value = 2366.45
$ws.Range("K55").Value = 1185.4286
$ws.Range("L55").Value = 2366.45
$ws.Range("M55").Value = -1012.4286
$ws.Range("N55").Value = -2712.45
$ws.Range("H61").Value = 6334.476
$ws.Range("I61").Value = 5546.625
$ws.Range("K61").Value = 5546.625
$ws.Range("M61").Value = -5344.625
$ws.Range("H68").Value = 2249.077
$ws.Range("I68").Value = 2019.8334
$ws.Range("K68").Value = 2019.8334
$ws.Range("M68").Value = -1270.8334
$ws.Range("H71").Value = 2249.077
$ws.Range("I71").Value = 2019.8334
$ws.Range("K71").Value = 10099.167
$ws.Range("M71").Value = -6355.166999999999
$ws.Range("H93").Value = 820.4583
$ws.Range("I93").Value = 1027.5714
$ws.Range("J93").Value = 530.5
$ws.Range("K93").Value = 1027.5714
$ws.Range("L93").Value = 530.5
$ws.Range("M93").Value = 220.4286
$ws.Range("N93").Value = -3026.5
$ws.Range("H100").Value = 2964.4583
$ws.Range("I100").Value = 6670.8
$ws.Range("K100").Value = 6670.8
$ws.Range("M100").Value = -6129.8
$ws.Range("H113").Value = 6334.476
$ws.Range("I113").Value = 5546.625
$ws.Range("K113").Value = 5546.625
$ws.Range("M113").Value = -3376.625
$ws.Range("H122").Value = 5651.143
$ws.Range("I122").Value = 3247.8
$ws.Range("J122").Value = 6986.3335
$ws.Range("K122").Value = 9743.400000000001
$ws.Range("L122").Value = 20959.0005
$ws.Range("M122").Value = -7293.400000000001
$ws.Range("N122").Value = -25859.0005
$ws.Range("H126").Value = 4616.931
$ws.Range("I126").Value = 3173.75
$ws.Range("J126").Value = 5635.647
$ws.Range("K126").Value = 9521.25
$ws.Range("L126").Value = 16906.941
$ws.Range("M126").Value = -7051.25
$ws.Range("N126").Value = -21846.941
$ws.Range("H132").Value = 2292.875
$ws.Range("I132").Value = 2292.875
$ws.Range("K132").Value = 6878.625
$ws.Range("M132").Value = -4348.625
$ws.Range("H136").Value = 13161603
$ws.Range("I136").Value = 62501588
$ws.Range("J136").Value = 4273.467
$ws.Range("K136").Value = 187504764
$ws.Range("L136").Value = 12820.401
$ws.Range("M136").Value = -187502214
$ws.Range("N136").Value = -17920.401
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("N27").ClearContents()
$ws.Range("N137").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 10406.25
$ws.Range("J45").Value = 10208.667
$ws.Range("L45").Value = 10208.667
$ws.Range("N45").Value = -11190.667
$ws.Range("H81").Value = 755.8889
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("H84").Value = 755.8889
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("H96").Value = 3199.875
$ws.Range("I96").Value = 2116.5
$ws.Range("J96").Value = 3354.6428
$ws.Range("K96").Value = 2116.5
$ws.Range("L96").Value = 3354.6428
$ws.Range("M96").Value = -743.5
$ws.Range("N96").Value = -6100.6428
$ws.Range("H100").Value = 1987.25
$ws.Range("I100").Value = 1483
$ws.Range("J100").Value = 3500
$ws.Range("K100").Value = 2966
$ws.Range("L100").Value = 7000
$ws.Range("M100").Value = -2425
$ws.Range("N100").Value = -8082
$ws.Range("H107").Value = 1114.5
$ws.Range("I107").Value = 819.6667
$ws.Range("J107").Value = 1999
$ws.Range("K107").Value = 2459.0001
$ws.Range("L107").Value = 5997
$ws.Range("M107").Value = -539.0001000000002
$ws.Range("N107").Value = -9837
$ws.Range("H122").Value = 2867.8235
$ws.Range("I122").Value = 2450.625
$ws.Range("J122").Value = 3238.6667
$ws.Range("K122").Value = 7351.875
$ws.Range("L122").Value = 9716.000100000001
$ws.Range("M122").Value = -4901.875
$ws.Range("N122").Value = -14616.0001
$ws.Range("H126").Value = 3990.158
$ws.Range("I126").Value = 2144.889
$ws.Range("J126").Value = 5650.9
$ws.Range("K126").Value = 6434.667
$ws.Range("L126").Value = 16952.7
$ws.Range("M126").Value = -3964.667
$ws.Range("N126").Value = -21892.7
$ws.Range("H132").Value = 5562.92
$ws.Range("I132").Value = 2007
$ws.Range("J132").Value = 10088.637
$ws.Range("K132").Value = 6021
$ws.Range("L132").Value = 30265.911
$ws.Range("M132").Value = -3491
$ws.Range("N132").Value = -35325.911
$ws.Range("H136").Value = 17898604
$ws.Range("I136").Value = 22774326
$ws.Range("J136").Value = 20956
$ws.Range("K136").Value = 68322978
$ws.Range("L136").Value = 62868
$ws.Range("M136").Value = -68320428
$ws.Range("N136").Value = -67968
$ws.Range("H139").Value = 117484
$ws.Range("J139").Value = 117484
$ws.Range("L139").Value = 117484
$ws.Range("N139").Value = -127764
$ws.Range("N81").ClearContents()
$ws.Range("N84").ClearContents()
